# Updates During The Period Till 05-04-2025
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Rename the second sheet
$ws2.Name = "SCAN accmu Real UseCase"

# Update the selection on the first sheet (no longer the active tab)
$ws1.Range("C7").Select() | Out-Null

# Make the renamed sheet the active / selected tab
$ws2.Activate() | Out-Null
